# Clean-up of input tables:
#  - rename the worksheet from "updated" to "Tabelle1"
#  - move the active cell selection from B5 to A5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("updated")

$ws.Name = "Tabelle1"

$ws.Activate()
$ws.Range("A5").Select()
